# Se agregan nuevos casos de uso:
# 1. Cesion de contrato nit a nit con cambio a plan pospago empresarial 5.3
# 2. Activacion nintendo con cliente nit
#
# This updates existing "caso de uso" rows (9-14) with new transaction
# references, adds a client note to rows 12/13, fixes the portability WSDL
# URL in row 2, and appends two brand-new rows (15/16) for the new cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 11")

# --- Row 2: portability SoapUI endpoint moved host/port and now asks for ?wsdl ---
$ws.Range("I2").Value = "http://10.65.50.8:8080/PortabilidadServiceEAR-HPNPCommunicationsDelegateEJB/NPCRMWSImpl?wsdl"

# --- Row 9: refreshed transaction / confirmation numbers ---
$ws.Range("C9").Value = "3052749177"
$ws.Range("D9").Value = "732111193280551"
$ws.Range("E9").Value = "3043209868"

# --- Row 10: refreshed transaction / confirmation numbers ---
$ws.Range("C10").Value = "3052754285"
$ws.Range("D10").Value = "732111324709512"
$ws.Range("E10").Value = "3046008593"

# --- Row 11: refreshed transaction / confirmation numbers (mirrors row 9) ---
$ws.Range("C11").Value = "3052749177"
$ws.Range("D11").Value = "732111193280551"

# --- Row 12: refreshed numbers + new "client nit a nit" annotation ---
$ws.Range("C12").Value = "3045987650"
$ws.Range("D12").Value = "732111324709673"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Font.Underline = $false
$ws.Range("E12").Value = "client nit a nit"

# --- Row 13: refreshed number + new annotation ---
$ws.Range("D13").Value = "732111324709674"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Font.Underline = $false
$ws.Range("E13").Value = "988154393"

# --- Row 14: refreshed number ---
$ws.Range("D14").Value = "732111324709675"

# --- Row 15 (new): "Cesion de contrato nit a nit" / plan pospago empresarial 5.3 ---
$r15 = $ws.Range("A15:D15")
$r15.NumberFormat = "@"
$r15.Font.Underline = $false
$ws.Range("A15").Value = "10960370"
$ws.Range("B15").Value = "36844580"
$ws.Range("C15").Value = "3052754289"
$ws.Range("D15").Value = "732111324709676"

# --- Row 16 (new): "Activacion nintendo" with cliente nit ---
$r16 = $ws.Range("A16:D16")
$r16.NumberFormat = "@"
$r16.Font.Underline = $false
$ws.Range("A16").Value = "10960370"
$ws.Range("B16").Value = "914355426"
$ws.Range("C16").Value = "3046008586"
$ws.Range("D16").Value = "732111193278871"

# --- View state: Excel scrolled back to the top-left and the last selection
#     now sits on the newly-added row 16 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("F16").Select()

Write-Output "done"
